$d = $word.ActiveDocument

# 1. Update the date text
$d.Content.Find.Execute("February 9, 2021", $true, $false, $false, $false, $false,
                         $true, 1, $false, "March 5, 2021", 2)

# 2. Split the "Thank you" sentence and insert a new sentence before it.
#    The original run text begins ". Thank you for taking..." (note: a
#    non-breaking space follows the period in the source document).
$d.Content.Find.Execute("Thank you for taking the time to review my application", $true, $false, $false, $false, $false,
                         $true, 1, $false, "As a data analyst at your company, I would be able to apply my experience to provide data insights that help businesses make informed decisions and reduce risks. Thank you for taking the time to review my application", 2)
